# Refresh Kujata leve profit/price snapshot values across all job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5424.615
$ws.Range("I28").Value = 5480.6665
$ws.Range("J28").Value = 5298.5
$ws.Range("K28").Value = 5480.6665
$ws.Range("L28").Value = 5298.5
$ws.Range("M28").Value = -4995.6665
$ws.Range("N28").Value = -6268.5
$ws.Range("H40").Value = 888.2
$ws.Range("I40").Value = 867.75
$ws.Range("J40").Value = 970
$ws.Range("K40").Value = 867.75
$ws.Range("L40").Value = 970
$ws.Range("M40").Value = -692.75
$ws.Range("N40").Value = -1320
$ws.Range("I43").Value = 936
$ws.Range("J43").Value = 6947570
$ws.Range("K43").Value = 936
$ws.Range("L43").Value = 6947570
$ws.Range("M43").Value = -867
$ws.Range("N43").Value = -6947708
$ws.Range("H62").Value = 22227422
$ws.Range("I62").Value = 37042036
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 37042036
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -37041412
$ws.Range("N62").Value = -6748
$ws.Range("H65").Value = 22227422
$ws.Range("I65").Value = 37042036
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 185210180
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -185207060
$ws.Range("N65").Value = -33740
$ws.Range("H98").Value = 3142.8
$ws.Range("I98").Value = 3395.0454
$ws.Range("K98").Value = 3395.0454
$ws.Range("M98").Value = -1897.0454
$ws.Range("H122").Value = 3142.8
$ws.Range("I122").Value = 3395.0454
$ws.Range("K122").Value = 10185.1362
$ws.Range("M122").Value = -7735.136200000001
$ws.Range("H137").Value = 3006.5
$ws.Range("I137").Value = 2648
$ws.Range("J137").Value = 3169.4546
$ws.Range("K137").Value = 7944
$ws.Range("L137").Value = 9508.363799999999
$ws.Range("M137").Value = -5394
$ws.Range("N137").Value = -14608.3638
$ws.Range("H138").Value = 2537.375
$ws.Range("I138").Value = 1485.8
$ws.Range("J138").Value = 2753.4521
$ws.Range("K138").Value = 4457.4
$ws.Range("L138").Value = 8260.356299999999
$ws.Range("M138").Value = 682.6000000000004
$ws.Range("N138").Value = -18540.3563

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10820.534
$ws.Range("I32").Value = 8289.303
$ws.Range("J32").Value = 18414.227
$ws.Range("K32").Value = 8289.303
$ws.Range("L32").Value = 18414.227
$ws.Range("M32").Value = -8002.303
$ws.Range("N32").Value = -18988.227
$ws.Range("H45").Value = 1182.4117
$ws.Range("I45").Value = 954.63635
$ws.Range("K45").Value = 954.63635
$ws.Range("M45").Value = -577.63635
$ws.Range("H61").Value = 100001496
$ws.Range("I61").Value = 142858000
$ws.Range("K61").Value = 142858000
$ws.Range("M61").Value = -142857788
$ws.Range("H136").Value = 100001496
$ws.Range("I136").Value = 142858000
$ws.Range("K136").Value = 428574000
$ws.Range("M136").Value = -428571450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 201982580
$ws.Range("I105").Value = 252477730
$ws.Range("K105").Value = 252477730
$ws.Range("M105").Value = -252475983
$ws.Range("H107").Value = 879.1667
$ws.Range("I107").Value = 825
$ws.Range("K107").Value = 825
$ws.Range("M107").Value = 1095
$ws.Range("H134").Value = 4714.077
$ws.Range("I134").Value = 797.2174
$ws.Range("J134").Value = 34743.332
$ws.Range("K134").Value = 2391.6522
$ws.Range("L134").Value = 104229.996
$ws.Range("M134").Value = 143.3478
$ws.Range("N134").Value = -109299.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1569.8536
$ws.Range("I31").Value = 1459.1
$ws.Range("J31").Value = 6000
$ws.Range("K31").Value = 1459.1
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = -1164.1
$ws.Range("N31").Value = -6590
$ws.Range("H34").Value = 1569.8536
$ws.Range("I34").Value = 1459.1
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 1459.1
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -1257.1
$ws.Range("N34").Value = -6404
$ws.Range("H58").Value = 8036.25
$ws.Range("I58").Value = 1083.2858
$ws.Range("K58").Value = 1083.2858
$ws.Range("M58").Value = -880.2858000000001
$ws.Range("H99").Value = 2393998
$ws.Range("I99").Value = 4387429.5
$ws.Range("J99").Value = 1879.8
$ws.Range("K99").Value = 4387429.5
$ws.Range("L99").Value = 1879.8
$ws.Range("M99").Value = -4385931.5
$ws.Range("N99").Value = -4875.8
$ws.Range("H107").Value = 982.2381
$ws.Range("I107").Value = 630.5333000000001
$ws.Range("K107").Value = 630.5333000000001
$ws.Range("M107").Value = 1289.4667
$ws.Range("H126").Value = 2393998
$ws.Range("I126").Value = 4387429.5
$ws.Range("J126").Value = 1879.8
$ws.Range("K126").Value = 13162288.5
$ws.Range("L126").Value = 5639.4
$ws.Range("M126").Value = -13159818.5
$ws.Range("N126").Value = -10579.4
$ws.Range("H136").Value = 8036.25
$ws.Range("I136").Value = 1083.2858
$ws.Range("K136").Value = 3249.8574
$ws.Range("M136").Value = -699.8574000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 697.0204
$ws.Range("I113").Value = 626.5
$ws.Range("J113").Value = 745.65515
$ws.Range("K113").Value = 1879.5
$ws.Range("L113").Value = 2236.96545
$ws.Range("M113").Value = 290.5
$ws.Range("N113").Value = -6576.96545
$ws.Range("H137").Value = 27784354
$ws.Range("I137").Value = 62502330
$ws.Range("J137").Value = 9976.866
$ws.Range("K137").Value = 187506990
$ws.Range("L137").Value = 29930.598
$ws.Range("M137").Value = -187501890
$ws.Range("N137").Value = -40130.598

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 112502504
$ws.Range("I70").Value = 83336670
$ws.Range("K70").Value = 83336670
$ws.Range("M70").Value = -83336400
$ws.Range("H73").Value = 112502504
$ws.Range("I73").Value = 83336670
$ws.Range("K73").Value = 83336670
$ws.Range("M73").Value = -83335734
$ws.Range("H134").Value = 24243.9
$ws.Range("J134").Value = 24243.9
$ws.Range("L134").Value = 72731.70000000001
$ws.Range("N134").Value = -77801.70000000001
$ws.Range("H136").Value = 25050.111
$ws.Range("J136").Value = 25050.111
$ws.Range("L136").Value = 75150.333
$ws.Range("N136").Value = -80250.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2378.7334
$ws.Range("I7").Value = 2196.5
$ws.Range("J7").Value = 2587
$ws.Range("K7").Value = 2196.5
$ws.Range("L7").Value = 2587
$ws.Range("M7").Value = -2084.5
$ws.Range("N7").Value = -2811
$ws.Range("H16").Value = 827.5925999999999
$ws.Range("I16").Value = 827.5925999999999
$ws.Range("K16").Value = 827.5925999999999
$ws.Range("M16").Value = -657.5925999999999
$ws.Range("H22").Value = 1403.6666
$ws.Range("I22").Value = 1105.5
$ws.Range("K22").Value = 1105.5
$ws.Range("M22").Value = -810.5
$ws.Range("H27").Value = 1403.6666
$ws.Range("I27").Value = 1105.5
$ws.Range("K27").Value = 1105.5
$ws.Range("M27").Value = -998.5
$ws.Range("H61").Value = 1546.2307
$ws.Range("I61").Value = 1463.2727
$ws.Range("J61").Value = 2002.5
$ws.Range("K61").Value = 1463.2727
$ws.Range("L61").Value = 2002.5
$ws.Range("M61").Value = -1261.2727
$ws.Range("N61").Value = -2406.5
$ws.Range("H113").Value = 1546.2307
$ws.Range("I113").Value = 1463.2727
$ws.Range("J113").Value = 2002.5
$ws.Range("K113").Value = 1463.2727
$ws.Range("L113").Value = 2002.5
$ws.Range("M113").Value = 706.7273
$ws.Range("N113").Value = -6342.5
$ws.Range("H126").Value = 2378.7334
$ws.Range("I126").Value = 2196.5
$ws.Range("J126").Value = 2587
$ws.Range("K126").Value = 6589.5
$ws.Range("L126").Value = 7761
$ws.Range("M126").Value = -4119.5
$ws.Range("N126").Value = -12701

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 125559.125
$ws.Range("I107").Value = 598.6
$ws.Range("J107").Value = 333826.66
$ws.Range("K107").Value = 1795.8
$ws.Range("L107").Value = 1001479.98
$ws.Range("M107").Value = 124.1999999999998
$ws.Range("N107").Value = -1005319.98
$ws.Range("H136").Value = 1071.1428
$ws.Range("I136").Value = 850.5625
$ws.Range("J136").Value = 1777
$ws.Range("K136").Value = 2551.6875
$ws.Range("L136").Value = 5331
$ws.Range("M136").Value = -1.6875
$ws.Range("N136").Value = -10431
